$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.228.98"
$ws.Range("E2").Value = "  +3.34%  "

$ws.Range("D3").Value = "'1.814.41"
$ws.Range("E3").Value = "  +4.58%  "

$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").Value = "'329.29"
$ws.Range("E5").Value = "  +2.05%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("D7").Value = "'0.4438"
$ws.Range("E7").Value = "  +4.53%  "

$ws.Range("D8").Value = "'0.3724"
$ws.Range("E8").Value = "  +4.03%  "

$ws.Range("D9").Value = "'44.90"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").Value = "'0.07711"
$ws.Range("E10").Value = "  +5.15%  "

$ws.Range("D11").Value = "'1.125"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").Value = "'22.04"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("D14").Value = "'6.268"
$ws.Range("E14").Value = "  +3.52%  "

$ws.Range("D15").Value = "'7.547"
$ws.Range("E15").Value = "  +5.49%  "

$ws.Range("D16").Value = "'1.816.75"
$ws.Range("E16").Value = "  +4.42%  "

$ws.Range("D17").Value = "'93.28"
$ws.Range("E17").Value = "  +11.00%  "

$ws.Range("D18").Value = "'0.00001082"
$ws.Range("E18").Value = "  +2.62%  "

$ws.Range("D19").Value = "'0.06526"
$ws.Range("E19").Value = "  +9.74%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").Value = "'17.52"
$ws.Range("E21").Value = "  +5.05%  "

$ws.Range("D22").Value = "'6.229"
$ws.Range("E22").Value = "  +4.09%  "

$ws.Range("D23").Value = "'0.5335"
$ws.Range("E23").Value = "  -1.70%  "

$ws.Range("D24").Value = "'28.288.03"
$ws.Range("E24").Value = "  +3.49%  "

$ws.Range("D25").Value = "'11.70"
$ws.Range("E25").Value = "  +4.21%  "

$ws.Range("D26").Value = "'2.046"
$ws.Range("E26").Value = "  -14.71%  "

$ws.Range("D27").Value = "'20.64"
$ws.Range("E27").Value = "  +4.28%  "

$ws.Range("D28").Value = "'155.17"
$ws.Range("E28").Value = "  +4.13%  "

$ws.Range("D29").Value = "'2.022.01"
$ws.Range("E29").Value = "  +4.23%  "

$ws.Range("D30").Value = "'2.321"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("D31").Value = "'127.75"
$ws.Range("E31").Value = "  +1.75%  "

$ws.Range("E32").Value = "  -5.66%  "

$ws.Range("D33").Value = "'5.867"
$ws.Range("E33").Value = "  +6.27%  "

$ws.Range("D34").Value = "'0.09201"
$ws.Range("E34").Value = "  +2.49%  "

$ws.Range("D35").Value = "'3.678"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").Value = "'13.08"
$ws.Range("E36").Value = "  +7.04%  "

$ws.Range("D37").Value = "'0.02348"
$ws.Range("E37").Value = "  +4.12%  "

$ws.Range("D38").Value = "'0.2169"
$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("D39").Value = "'5.179"
$ws.Range("E39").Value = "  +4.46%  "

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.6569"
$ws.Range("E40").Value = "  +2.76%  "

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = "'0.06200"
$ws.Range("E41").Value = "  +2.15%  "

$ws.Range("D42").Value = "'1.200"
$ws.Range("E42").Value = "  +1.93%  "

$ws.Range("D43").Value = "'8.093"
$ws.Range("E43").Value = "  +3.66%  "

$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = "'1.390"
$ws.Range("E45").Value = "  -1.58%  "

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'13.85"
$ws.Range("E46").Value = "  +2.97%  "

$ws.Range("D47").Value = "'0.6075"
$ws.Range("E47").Value = "  +3.98%  "

$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'2.037"
$ws.Range("E49").Value = "  +5.68%  "

$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'126.58"
$ws.Range("E50").Value = "  +2.16%  "

$ws.Range("D51").Value = "'0.06987"
$ws.Range("E51").Value = "  +2.80%  "
